$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.076417455671822
$ws.Range("D2").Value = 1.077045107061167
$ws.Range("E2").Value = 1.079732698394253
$ws.Range("F2").Value = 1.08945147429767
$ws.Range("I2").Value = 1.062188075918561
$ws.Range("J2").Value = 1.081317015860643
$ws.Range("K2").Value = 1.079727440783993
$ws.Range("L2").Value = 1.082407978077283
$ws.Range("M2").Value = 1.092101558835843

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.077657705553621
$ws.Range("D3").Value = 1.078034300536449
$ws.Range("E3").Value = 1.080831446594123
$ws.Range("F3").Value = 1.090603890804163
$ws.Range("I3").Value = 1.062624509904485
$ws.Range("J3").Value = 1.08221572525808
$ws.Range("K3").Value = 1.080533833884042
$ws.Range("L3").Value = 1.083324169501241
$ws.Range("M3").Value = 1.093073114497412

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.0784599984479
$ws.Range("D4").Value = 1.078674096287257
$ws.Range("E4").Value = 1.081542359179416
$ws.Range("F4").Value = 1.091349655251975
$ws.Range("I4").Value = 1.062905487940874
$ws.Range("J4").Value = 1.082796452219712
$ws.Range("K4").Value = 1.081054726057717
$ws.Range("L4").Value = 1.083916357163479
$ws.Range("M4").Value = 1.093701259569581

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.078797228890585
$ws.Range("D5").Value = 1.078943000670841
$ws.Range("E5").Value = 1.081841215792141
$ws.Range("F5").Value = 1.09166319352433
$ws.Range("I5").Value = 1.063023270881851
$ws.Range("J5").Value = 1.083040399956981
$ws.Range("K5").Value = 1.081273495118513
$ws.Range("L5").Value = 1.084165158499748
$ws.Range("M5").Value = 1.093965209261527

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.07885384828678
$ws.Range("D6").Value = 1.078988147071811
$ws.Range("E6").Value = 1.081891394519036
$ws.Range("F6").Value = 1.091715839123143
$ws.Range("I6").Value = 1.063043027224346
$ws.Range("J6").Value = 1.083081348755668
$ws.Range("K6").Value = 1.081310214885764
$ws.Range("L6").Value = 1.084206924319999
$ws.Range("M6").Value = 1.094009520425924

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.078464504745276
$ws.Range("D7").Value = 1.078677689657971
$ws.Range("E7").Value = 1.081546352558018
$ws.Range("F7").Value = 1.09135384468995
$ws.Range("I7").Value = 1.062907063098317
$ws.Range("J7").Value = 1.082799712605917
$ws.Range("K7").Value = 1.081057650100284
$ws.Range("L7").Value = 1.083919682264226
$ws.Range("M7").Value = 1.093704786958196

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.076836653096697
$ws.Range("D8").Value = 1.077379467874875
$ws.Range("E8").Value = 1.080104036747776
$ws.Range("F8").Value = 1.089840924273873
$ws.Range("I8").Value = 1.062335866064666
$ws.Range("J8").Value = 1.081620904881902
$ws.Range("K8").Value = 1.080000151660504
$ws.Range("L8").Value = 1.082717744408411
$ws.Range("M8").Value = 1.092430008227801

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.073966306951663
$ws.Range("D9").Value = 1.075089659131926
$ws.Range("E9").Value = 1.077562049143758
$ws.Range("F9").Value = 1.087175476936277
$ws.Range("I9").Value = 1.061318408194648
$ws.Range("J9").Value = 1.079537528785603
$ws.Range("K9").Value = 1.0781297812344
$ws.Range("L9").Value = 1.080594751724219
$ws.Range("M9").Value = 1.090179679142215

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.072051363401263
$ws.Range("D10").Value = 1.073561594855487
$ws.Range("E10").Value = 1.075867017092999
$ws.Range("F10").Value = 1.08539877756378
$ws.Range("I10").Value = 1.060632707004902
$ws.Range("J10").Value = 1.078144388348381
$ws.Range("K10").Value = 1.076878148210886
$ws.Range("L10").Value = 1.079175972116246
$ws.Range("M10").Value = 1.088676692623241

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.071221815669955
$ws.Range("D11").Value = 1.072899549814288
$ws.Range("E11").Value = 1.075132943201026
$ws.Range("F11").Value = 1.084629491800849
$ws.Range("I11").Value = 1.060334026506331
$ws.Range("J11").Value = 1.077540123653932
$ws.Range("K11").Value = 1.076335042846464
$ws.Range("L11").Value = 1.078560789406904
$ws.Range("M11").Value = 1.088025209415536

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.070913627358375
$ws.Range("D12").Value = 1.072653577890865
$ws.Range("E12").Value = 1.074860256750283
$ws.Range("F12").Value = 1.084343748666463
$ws.Range("I12").Value = 1.060222816804728
$ws.Range("J12").Value = 1.07731551711052
$ws.Range("K12").Value = 1.076133136669698
$ws.Range("L12").Value = 1.078332155139995
$ws.Range("M12").Value = 1.087783115791312

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.070979737434005
$ws.Range("D13").Value = 1.072706342391393
$ws.Range("E13").Value = 1.074918749792933
$ws.Range("F13").Value = 1.084405041391497
$ws.Range("I13").Value = 1.060246683736469
$ws.Range("J13").Value = 1.077363703049927
$ws.Range("K13").Value = 1.076176454065693
$ws.Range("L13").Value = 1.078381203784353
$ws.Range("M13").Value = 1.087835050415863

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.071196341912811
$ws.Range("D14").Value = 1.07287921890813
$ws.Range("E14").Value = 1.075110403245265
$ws.Range("F14").Value = 1.084605872115386
$ws.Range("I14").Value = 1.060324839318962
$ws.Range("J14").Value = 1.077521560793303
$ws.Range("K14").Value = 1.076318356749219
$ws.Range("L14").Value = 1.078541893047573
$ws.Range("M14").Value = 1.088005200006681

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.071329791404556
$ws.Range("D15").Value = 1.072985725987923
$ws.Range("E15").Value = 1.075228484733133
$ws.Range("F15").Value = 1.084729611012623
$ws.Range("I15").Value = 1.060372958215378
$ws.Range("J15").Value = 1.077618801474216
$ws.Range("K15").Value = 1.076405764754832
$ws.Range("L15").Value = 1.078640881996068
$ws.Range("M15").Value = 1.088110020989794

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.072106409851537
$ws.Range("D16").Value = 1.073605524416701
$ws.Range("E16").Value = 1.07591573262306
$ws.Range("F16").Value = 1.085449833158362
$ws.Range("I16").Value = 1.060652492129476
$ws.Range("J16").Value = 1.078184469656635
$ws.Range("K16").Value = 1.076914168188079
$ws.Range("L16").Value = 1.079216781899078
$ws.Range("M16").Value = 1.088719914954508

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.07259346266666
$ws.Range("D17").Value = 1.073994204065062
$ws.Range("E17").Value = 1.076346793448365
$ws.Range("F17").Value = 1.085901618148782
$ws.Range("I17").Value = 1.060827362532792
$ws.Range("J17").Value = 1.078539022603782
$ws.Range("K17").Value = 1.077232770037238
$ws.Range("L17").Value = 1.079577802231444
$ws.Range("M17").Value = 1.089102302126545

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.072877517336562
$ws.Range("D18").Value = 1.07422087745496
$ws.Range("E18").Value = 1.076598213062869
$ws.Range("F18").Value = 1.086165140354479
$ws.Range("I18").Value = 1.060929191020152
$ws.Range("J18").Value = 1.078745728486331
$ws.Range("K18").Value = 1.077418495206477
$ws.Range("L18").Value = 1.079788298169971
$ws.Range("M18").Value = 1.089325276579592

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.072974366779837
$ws.Range("D19").Value = 1.074298161009386
$ws.Range("E19").Value = 1.076683938858601
$ws.Range("F19").Value = 1.086254995288318
$ws.Range("I19").Value = 1.060963883018011
$ws.Range("J19").Value = 1.078816193132276
$ws.Range("K19").Value = 1.077481804090449
$ws.Range("L19").Value = 1.079860058133519
$ws.Range("M19").Value = 1.089401294040516

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.072541210148803
$ws.Range("D20").Value = 1.073952506241921
$ws.Range("E20").Value = 1.07630054583861
$ws.Range("F20").Value = 1.085853145529121
$ws.Range("I20").Value = 1.060808618232338
$ws.Range("J20").Value = 1.07850099266557
$ws.Range("K20").Value = 1.077198598452311
$ws.Range("L20").Value = 1.079539076568515
$ws.Range("M20").Value = 1.089061282390707

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.071132558896373
$ws.Range("D21").Value = 1.072828312720118
$ws.Range("E21").Value = 1.075053966610263
$ws.Range("F21").Value = 1.084546732379971
$ws.Range("I21").Value = 1.060301831811205
$ws.Range("J21").Value = 1.077475079934442
$ws.Range("K21").Value = 1.076276574723261
$ws.Range("L21").Value = 1.078494577611006
$ws.Range("M21").Value = 1.087955098069801

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.070246550863372
$ws.Range("D22").Value = 1.072121146080094
$ws.Range("E22").Value = 1.074270083644394
$ws.Range("F22").Value = 1.083725360026065
$ws.Range("I22").Value = 1.059981652919771
$ws.Range("J22").Value = 1.076829146060217
$ws.Range("K22").Value = 1.075695862147033
$ws.Range("L22").Value = 1.077837118220495
$ws.Range("M22").Value = 1.087258995016049

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.070716272920762
$ws.Range("D23").Value = 1.072496061325797
$ws.Range("E23").Value = 1.074685645722697
$ws.Range("F23").Value = 1.084160783580201
$ws.Range("I23").Value = 1.060151532300699
$ws.Range("J23").Value = 1.077171653856723
$ws.Range("K23").Value = 1.076003804092504
$ws.Range("L23").Value = 1.078185720681348
$ws.Range("M23").Value = 1.08762806991572

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.072564820920592
$ws.Range("D24").Value = 1.073971347807368
$ws.Range("E24").Value = 1.076321443176439
$ws.Range("F24").Value = 1.085875048208242
$ws.Range("I24").Value = 1.060817088501598
$ws.Range("J24").Value = 1.078518177063161
$ws.Range("K24").Value = 1.077214039460271
$ws.Range("L24").Value = 1.079556575279358
$ws.Range("M24").Value = 1.089079817646162

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.074708594091094
$ws.Range("D25").Value = 1.075681893830621
$ws.Range("E25").Value = 1.07821927415668
$ws.Range("F25").Value = 1.087864506396157
$ws.Range("I25").Value = 1.061582745846135
$ws.Range("J25").Value = 1.080076870119242
$ws.Range("K25").Value = 1.078614143897532
$ws.Range("L25").Value = 1.081144199161515
$ws.Range("M25").Value = 1.090761925468654
